# Serramazzoni.xlsx update — "aggiornato a 2/3, aggiornati i report"
#
# A new daily data point (2021-02-08 / Excel serial 44235) was inserted into
# the time series, shifting every following row down by one. The rolling
# 7-day-sum columns (C = somma mobile 7gg., D = per-100k-abitanti) were
# recalculated for the rows whose trailing 7-day window now includes the new
# point, and two more days of fresh data (44256, 44257) were appended at the
# bottom (their C/D cells stay blank until enough future days accumulate,
# matching the source data's existing "blank" rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recalculated rolling-sum rows that keep their current date (90-92) ---
# row -> (B, C, D)
$recalc = @{
    90 = @(1, 14, 164.9581713208436)
    91 = @(1, 12, 141.3927182750088)
    92 = @(1, 12, 141.3927182750088)
}
foreach ($r in $recalc.Keys) {
    $vals = $recalc[$r]
    $ws.Cells.Item($r, 2).Value2 = $vals[0]
    $ws.Cells.Item($r, 3).Value2 = $vals[1]
    $ws.Cells.Item($r, 4).Value2 = $vals[2]
}

# --- Shift existing rows 93..113 down to 94..114 (bottom-up so we never
#     clobber a source row before it has been copied) ---
for ($r = 113; $r -ge 93; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dest, 4).Value2 = $ws.Cells.Item($r, 4).Value2
}

# --- New row 93: the inserted date, 2021-02-08 (serial 44235) ---
# (row 93 already exists/has the date style from before the shift, so a
# plain value write keeps its existing "s=2" cell style)
$ws.Cells.Item(93, 1).Value2 = 44235
$ws.Cells.Item(93, 2).Value2 = 1
$ws.Cells.Item(93, 3).Value2 = 6
$ws.Cells.Item(93, 4).Value2 = 70.69635913750442

# --- Row 112 (previously blank C/D, now enough data exists to fill it in) ---
$ws.Cells.Item(112, 3).Value2 = 25
$ws.Cells.Item(112, 4).Value2 = 294.5681630729351

# --- Append two brand-new trailing rows: 114 (44256) and 115 (44257) ---
$ws.Range("A113").Copy() | Out-Null
$ws.Range("A114").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(114, 1).Value2 = 44256
$ws.Cells.Item(114, 2).Value2 = 4

$ws.Range("A113").Copy() | Out-Null
$ws.Range("A115").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(115, 1).Value2 = 44257
$ws.Cells.Item(115, 2).Value2 = 10
